$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 9
    $ws.Cells.Item(9, 7).Value = 2.22   # G9: 2.25 -> 2.22
    $ws.Cells.Item(9, 8).Value = 2.75   # H9: 2.7 -> 2.75
    $ws.Cells.Item(9, 9).Value = 3.7   # I9: 3.75 -> 3.7
    $ws.Cells.Item(9, 10).Value = 1.09   # J9: 1.1 -> 1.09
    $ws.Cells.Item(9, 11).Value = 6.2   # K9: 5.9 -> 6.2
    $ws.Cells.Item(9, 12).Value = 1.35   # L9: 1.38 -> 1.35
    $ws.Cells.Item(9, 13).Value = 2.95   # M9: 2.82 -> 2.95
    $ws.Cells.Item(9, 14).Value = 2.05   # N9: 2.12 -> 2.05
    $ws.Cells.Item(9, 15).Value = 1.7   # O9: 1.65 -> 1.7
    $ws.Cells.Item(9, 16).Value = 1.42   # P9: 1.44 -> 1.42
    $ws.Cells.Item(9, 17).Value = 2.67   # Q9: 2.62 -> 2.67
    $ws.Cells.Item(9, 18).Value = 1.7   # R9: 1.75 -> 1.7
    $ws.Cells.Item(9, 19).Value = 2.02   # S9: 1.98 -> 2.02
    $ws.Cells.Item(9, 20).Value = 7.2   # T9: 6.9 -> 7.2
    $ws.Cells.Item(9, 21).Value = 11.25   # U9: 11 -> 11.25
    $ws.Cells.Item(9, 23).Value = 24   # W9: 25 -> 24
    $ws.Cells.Item(9, 24).Value = 18.5   # X9: 19 -> 18.5
    $ws.Cells.Item(9, 25).Value = 26   # Y9: 28 -> 26
    $ws.Cells.Item(9, 26).Value = 6.2   # Z9: 5.9 -> 6.2
    $ws.Cells.Item(9, 27).Value = 5.4   # AA9: 5.3 -> 5.4
    $ws.Cells.Item(9, 28).Value = 12   # AB9: 12.5 -> 12
    $ws.Cells.Item(9, 29).Value = 55   # AC9: 60 -> 55
    $ws.Cells.Item(9, 30).Value = 450   # AD9: 500 -> 450
    $ws.Cells.Item(9, 31).Value = 10   # AE9: 9.75 -> 10
    $ws.Cells.Item(9, 33).Value = 11.75   # AG9: 12 -> 11.75
    $ws.Cells.Item(9, 34).Value = 60   # AH9: 65 -> 60
    $ws.Cells.Item(9, 35).Value = 35   # AI9: 37 -> 35
    $ws.Cells.Item(9, 36).Value = 37   # AJ9: 40 -> 37
    # Row 12
    $ws.Cells.Item(12, 7).Value = 3.3   # G12: 3.5 -> 3.3
    $ws.Cells.Item(12, 9).Value = 2.05   # I12: 1.95 -> 2.05
    $ws.Cells.Item(12, 11).Value = 9.5   # K12: 10 -> 9.5
    $ws.Cells.Item(12, 28).Value = 17   # AB12: 19 -> 17
    # Row 13
    $ws.Cells.Item(13, 11).Value = 12   # K13: 13 -> 12
    $ws.Cells.Item(13, 20).Value = 7.5   # T13: 8 -> 7.5
    $ws.Cells.Item(13, 21).Value = 8.5   # U13: 9 -> 8.5
    $ws.Cells.Item(13, 29).Value = 51   # AC13: 41 -> 51
    $ws.Cells.Item(13, 36).Value = 41   # AJ13: 34 -> 41
    # Row 14
    $ws.Cells.Item(14, 14).Value = 1.75   # N14: 1.73 -> 1.75
    $ws.Cells.Item(14, 15).Value = 2.05   # O14: 2.08 -> 2.05
    # Row 21
    $ws.Cells.Item(21, 7).Value = 2.75   # G21: 2.8 -> 2.75
    $ws.Cells.Item(21, 22).Value = 11   # V21: 12 -> 11
    # Row 23
    $ws.Cells.Item(23, 8).Value = 3.5   # H23: 3.35 -> 3.5
    $ws.Cells.Item(23, 9).Value = 2.05   # I23: 2.07 -> 2.05
    $ws.Cells.Item(23, 14).Value = 1.7   # N23: 1.8 -> 1.7
    $ws.Cells.Item(23, 15).Value = 1.91   # O23: 1.8 -> 1.91
    $ws.Cells.Item(23, 20).Value = 9.25   # T23: 8.5 -> 9.25
    $ws.Cells.Item(23, 21).Value = 14   # U23: 13 -> 14
    $ws.Cells.Item(23, 24).Value = 19.5   # X23: 20 -> 19.5
    $ws.Cells.Item(23, 25).Value = 23   # Y23: 25 -> 23
    $ws.Cells.Item(23, 26).Value = 11.75   # Z23: 10.25 -> 11.75
    $ws.Cells.Item(23, 27).Value = 6   # AA23: 5.8 -> 6
    $ws.Cells.Item(23, 28).Value = 11   # AB23: 11.5 -> 11
    $ws.Cells.Item(23, 29).Value = 40   # AC23: 45 -> 40
    $ws.Cells.Item(23, 30).Value = 250   # AD23: 300 -> 250
    $ws.Cells.Item(23, 31).Value = 7.3   # AE23: 6.7 -> 7.3
    $ws.Cells.Item(23, 33).Value = 7.4   # AG23: 7.5 -> 7.4
    $ws.Cells.Item(23, 35).Value = 13   # AI23: 13.5 -> 13
    $ws.Cells.Item(23, 36).Value = 19   # AJ23: 21 -> 19
    # Row 26
    $ws.Cells.Item(26, 7).Value = 3.25   # G26: 3.1 -> 3.25
    $ws.Cells.Item(26, 9).Value = 2.25   # I26: 2.35 -> 2.25
    $ws.Cells.Item(26, 11).Value = 8.5   # K26: 9 -> 8.5
    $ws.Cells.Item(26, 20).Value = 9.5   # T26: 9 -> 9.5
    $ws.Cells.Item(26, 24).Value = 29   # X26: 26 -> 29
    $ws.Cells.Item(26, 25).Value = 41   # Y26: 34 -> 41
    $ws.Cells.Item(26, 35).Value = 19   # AI26: 21 -> 19
    # Row 33
    $ws.Cells.Item(33, 10).Value = 1.01   # J33: 26 -> 1.01
    $ws.Cells.Item(33, 11).Value = 13   # K33: 1.02 -> 13
    # Row 39
    $ws.Cells.Item(39, 11).Value = 10   # K39: 9.5 -> 10
    $ws.Cells.Item(39, 14).Value = 2.08   # N39: 2.05 -> 2.08
    $ws.Cells.Item(39, 15).Value = 1.73   # O39: 1.75 -> 1.73
    # Row 40
    $ws.Cells.Item(40, 7).Value = 3.2   # G40: 3.1 -> 3.2
    $ws.Cells.Item(40, 8).Value = 3.6   # H40: 3.7 -> 3.6
    $ws.Cells.Item(40, 12).Value = 1.2   # L40: 1.22 -> 1.2
    $ws.Cells.Item(40, 13).Value = 4.33   # M40: 4 -> 4.33
    $ws.Cells.Item(40, 25).Value = 29   # Y40: 26 -> 29
    $ws.Cells.Item(40, 32).Value = 11   # AF40: 12 -> 11
    $ws.Cells.Item(40, 35).Value = 15   # AI40: 17 -> 15
    # Row 42
    $ws.Cells.Item(42, 7).Value = 1.95   # G42: 2.05 -> 1.95
    $ws.Cells.Item(42, 8).Value = 3.7   # H42: 3.6 -> 3.7
    $ws.Cells.Item(42, 9).Value = 3.5   # I42: 3.3 -> 3.5
    $ws.Cells.Item(42, 14).Value = 1.73   # N42: 1.8 -> 1.73
    $ws.Cells.Item(42, 15).Value = 2.08   # O42: 2 -> 2.08
    $ws.Cells.Item(42, 22).Value = 8.5   # V42: 9 -> 8.5
    $ws.Cells.Item(42, 23).Value = 17   # W42: 19 -> 17
    $ws.Cells.Item(42, 26).Value = 13   # Z42: 12 -> 13
